$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "26.141.13"
Set-TextValue "E2" "  -1.40%  "

# Row 3
Set-TextValue "D3" "1.651.81"
Set-TextValue "E3" "  -1.77%  "

# Row 4
Set-TextValue "E4" "  +0.35%  "

# Row 5
Set-TextValue "D5" "218.21"
Set-TextValue "E5" "  +0.38%  "

# Row 6
Set-TextValue "D6" "0.5208"
Set-TextValue "E6" "  -2.21%  "

# Row 7
Set-TextValue "E7" "  +0.31%  "

# Row 8
Set-TextValue "D8" "0.2662"
Set-TextValue "E8" "  -0.49%  "

# Row 9
Set-TextValue "E9" "  -1.67%  "

# Row 10
Set-TextValue "D10" "21.08"
Set-TextValue "E10" "  -1.84%  "

# Row 11
Set-TextValue "D11" "0.07713"
Set-TextValue "E11" "  -0.99%  "

# Row 12
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.656.81"
Set-TextValue "E12" "  -1.59%  "

# Row 13
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.425"
Set-TextValue "E13" "  -1.85%  "

# Row 14
Set-TextValue "D14" "1.879.71"

# Row 15
Set-TextValue "D15" "0.5464"
Set-TextValue "E15" "  -2.79%  "

# Row 16
Set-TextValue "D16" "0.0₅8239"
Set-TextValue "E16" "  -2.18%  "

# Row 17
Set-TextValue "D17" "64.82"
Set-TextValue "E17" "  -1.73%  "

# Row 18
Set-TextValue "D18" "26.168.62"
Set-TextValue "E18" "  -1.43%  "

# Row 19
Set-TextValue "D19" "1.005"
Set-TextValue "E19" "  +0.29%  "

# Row 20
Set-TextValue "D20" "4.660"
Set-TextValue "E20" "  -2.93%  "

# Row 21
Set-TextValue "D21" "192.95"
Set-TextValue "E21" "  -1.24%  "

# Row 22
Set-TextValue "D22" "10.16"
Set-TextValue "E22" "  -2.43%  "

# Row 23
Set-TextValue "D23" "6.089"
Set-TextValue "E23" "  -4.60%  "

# Row 24
Set-TextValue "D24" "1.007"

# Row 25
Set-TextValue "D25" "137.32"
Set-TextValue "E25" "  -4.12%  "

# Row 26
Set-TextValue "D26" "0.1242"
Set-TextValue "E26" "  -2.84%  "

# Row 27
Set-TextValue "D27" "7.230"
Set-TextValue "E27" "  -3.21%  "

# Row 28
Set-TextValue "D28" "16.12"
Set-TextValue "E28" "  -0.36%  "

# Row 29
Set-TextValue "D29" "1.430"
Set-TextValue "E29" "  +1.03%  "

# Row 30
Set-TextValue "D30" "0.06024"
Set-TextValue "E30" "  -1.61%  "

# Row 31
Set-TextValue "D31" "1.281"
Set-TextValue "E31" "  +0.25%  "

# Row 32
Set-TextValue "D32" "3.559"
Set-TextValue "E32" "  -1.30%  "

# Row 33
Set-TextValue "D33" "3.329"
Set-TextValue "E33" "  -3.72%  "

# Row 34
Set-TextValue "D34" "1.646"
Set-TextValue "E34" "  -3.38%  "

# Row 35
Set-TextValue "D35" "0.9800"
Set-TextValue "E35" "  -3.46%  "

# Row 36
Set-TextValue "D36" "2.410"
Set-TextValue "E36" "  -0.35%  "

# Row 37
Set-TextValue "D37" "2.768"
Set-TextValue "E37" "  -0.86%  "

# Row 38
Set-TextValue "D38" "0.5917"
Set-TextValue "E38" "  +3.83%  "

# Row 39
Set-TextValue "D39" "0.01590"
Set-TextValue "E39" "  -3.15%  "

# Row 40
Set-TextValue "D40" "5.954"
Set-TextValue "E40" "  +0.03%  "

# Row 41
Set-TextValue "D41" "0.8617"
Set-TextValue "E41" "  -1.02%  "

# Row 42
Set-TextValue "D42" "1.003"
Set-TextValue "E42" "  +0.14%  "

# Row 43
Set-TextValue "D43" "1.041.45"
Set-TextValue "E43" "  -1.85%  "

# Row 44
Set-TextValue "D44" "99.55"
Set-TextValue "E44" "  -0.39%  "

# Row 45
Set-TextValue "D45" "1.793.61"
Set-TextValue "E45" "  -2.09%  "

# Row 46
Set-TextValue "D46" "0.0₈112"
Set-TextValue "E46" "  +0.46%  "

# Row 47
Set-TextValue "D47" "57.08"
Set-TextValue "E47" "  -0.25%  "

# Row 48
Set-TextValue "D48" "1.005"
Set-TextValue "E48" "  +0.70%  "

# Row 49
Set-TextValue "D49" "8.115"
Set-TextValue "E49" "  -0.52%  "

# Row 50
Set-TextValue "D50" "0.05180"
Set-TextValue "E50" "  -0.42%  "

# Row 51
Set-TextValue "D51" "1.465"
Set-TextValue "E51" "  +4.16%  "
